# Update cryptocurrency price/volume data scraped on Wed Jan 31 09:56:42 UTC 2024
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range("D2").Value = '42.812.61'
$ws.Range("E2").Value = '  -1.30%  '

# Row 3: Ethereum
$ws.Range("D3").Value = '2.323.27'
$ws.Range("E3").Value = '  +0.75%  '

# Row 4: TetherUSD
$ws.Range("E4").Value = '  +0.04%  '

# Row 5: BNB
$ws.Range("D5").Value = '''304.71'
$ws.Range("E5").Value = '  -1.94%  '

# Row 6: Solana
$ws.Range("D6").Value = '''100.05'
$ws.Range("E6").Value = '  -3.52%  '

# Row 7: XRP
$ws.Range("E7").Value = '  -4.69%  '

# Row 8: USDC
$ws.Range("E8").Value = '  +0.05%  '

# Row 9: Cardano
$ws.Range("E9").Value = '  -4.40%  '

# Row 10: Avalanche
$ws.Range("D10").Value = '''34.49'
$ws.Range("E10").Value = '  -5.82%  '

# Row 11: OKB
$ws.Range("D11").Value = '''52.16'
$ws.Range("E11").Value = '  -0.72%  '

# Row 12: Dogecoin
$ws.Range("D12").Value = '''0.0793'
$ws.Range("E12").Value = '  -2.28%  '

# Row 13: TRON
$ws.Range("E13").Value = '  +1.02%  '

# Row 14: Polkadot
$ws.Range("D14").Value = '''6.77'
$ws.Range("E14").Value = '  -3.46%  '

# Row 15: WrappedliquidstakedEther2.0
$ws.Range("D15").Value = '2.689.17'
$ws.Range("E15").Value = '  +0.94%  '

# Row 16: Chainlink
$ws.Range("D16").Value = '''15.72'
$ws.Range("E16").Value = '  +4.21%  '

# Row 17: WrappedEther
$ws.Range("D17").Value = '2.297.20'
$ws.Range("E17").Value = '  -0.97%  '

# Row 18: Polygon
$ws.Range("D18").Value = '''0.826'
$ws.Range("E18").Value = '  +2.01%  '

# Row 19: WrappedBTC
$ws.Range("D19").Value = '42.759.40'
$ws.Range("E19").Value = '  -1.17%  '

# Row 21: Uniswap
$ws.Range("E21").Value = '  -0.36%  '

# Row 22: InternetComputer(DFINITY)
$ws.Range("D22").Value = '''11.60'
$ws.Range("E22").Value = '  -4.60%  '

# Row 23: Litecoin
$ws.Range("D23").Value = '''69.40'
$ws.Range("E23").Value = '  +1.94%  '

# Row 24: BitcoinCash
$ws.Range("D24").Value = '''235.48'
$ws.Range("E24").Value = '  -2.96%  '

# Row 25: ImmutableX
$ws.Range("E25").Value = '  -1.71%  '

# Row 26: PancakeSwap
$ws.Range("E26").Value = '  -3.27%  '

# Row 27: Dai
$ws.Range("E27").Value = '  +0.11%  '

# Row 28: EthereumClassic
$ws.Range("D28").Value = '''25.54'
$ws.Range("E28").Value = '  +2.66%  '

# Row 29: LEO
$ws.Range("E29").Value = '  -0.30%  '

# Row 30: Toncoin
$ws.Range("D30").Value = '''2.30'
$ws.Range("E30").Value = '  -0.06%  '

# Row 31: InjectiveProtocol
$ws.Range("D31").Value = '''34.70'
$ws.Range("E31").Value = '  -5.73%  '

# Row 32: Cosmos
$ws.Range("E32").Value = '  -4.38%  '

# Row 33: Monero
$ws.Range("D33").Value = '''160.07'

# Row 34: FirstDigitalUSD
$ws.Range("E34").Value = '  +0.02%  '

# Row 35: Filecoin
$ws.Range("E35").Value = '  -4.20%  '

# Row 36: RenderToken
$ws.Range("D36").Value = '''4.61'
$ws.Range("E36").Value = '  +3.64%  '

# Row 37: WEMIXToken
$ws.Range("E37").Value = '  -3.35%  '

# Row 38: Hedera
$ws.Range("D38").Value = '''0.0719'
$ws.Range("E38").Value = '  -3.10%  '

# Row 39: Celestia
$ws.Range("D39").Value = '''16.97'
$ws.Range("E39").Value = '  -6.73%  '

# Row 40: LidoDAOToken
$ws.Range("D40").Value = '''2.89'
$ws.Range("E40").Value = '  -5.40%  '

# Row 41: ARBITRUM
$ws.Range("E41").Value = '  -2.80%  '

# Row 42: Kaspa
$ws.Range("E42").Value = '  -4.62%  '

# Row 43: Stellar
$ws.Range("E43").Value = '  -3.16%  '

# Row 44: ApeXProtocol
$ws.Range("E44").Value = '  -4.16%  '

# Row 45: Maker
$ws.Range("D45").Value = '2.003.92'
$ws.Range("E45").Value = '  +1.12%  '

# Row 46: VeChain
$ws.Range("E46").Value = '  -3.87%  '

# Row 47: EnergySwap
$ws.Range("D47").Value = '''18.74'
$ws.Range("E47").Value = '  -1.62%  '

# Row 48: FraxShare
$ws.Range("D48").Value = '''10.18'
$ws.Range("E48").Value = '  +1.93%  '

# Row 49: NEARProtocol
$ws.Range("E49").Value = '  -4.58%  '

# Row 50: MultiversX
$ws.Range("D50").Value = '''55.55'
$ws.Range("E50").Value = '  -0.62%  '

# Row 51: HuobiToken
$ws.Range("D51").Value = '''2.89'
$ws.Range("E51").Value = '  -1.88%  '
